$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the new cells in the same order the source workbook's shared-string
# table implies (ItemName/Computers, daysOut/2, then the price/qty/UoM block).
$ws.Range("C1").Value = "ItemName"
$ws.Range("C2").Value = "Computers"
$ws.Range("G1").Value = "daysOut"
$ws.Range("G2").Value = "'2"
$ws.Range("D1").Value = "UnitPrice"
$ws.Range("E1").Value = "Quantity"
$ws.Range("F1").Value = "UnitofMeasure"
$ws.Range("F2").Value = "EA-EACH"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1

# Column widths (closest achievable values through the ColumnWidth API -
# the host quantizes stored width to 1/6-character steps, so these land on
# the nearest reachable value to the authored 10.7109375 / 14.42578125 / 8.28515625)
$ws.Columns.Item(3).ColumnWidth = 9.833333333333334
$ws.Columns.Item(4).ColumnWidth = 9.833333333333334
$ws.Columns.Item(5).ColumnWidth = 9.833333333333334
$ws.Columns.Item(6).ColumnWidth = 13.666666666666666
$ws.Columns.Item(7).ColumnWidth = 7.5

# Selection matches the authored file
$ws.Range("F7").Select()
